# Updated via Streamlit Approval System
# Sets APPROVAL_1 (AI) / APPROVAL_2 (AJ) decisions for each pending row,
# and zeroes out the COST_CENTER/LEDGER_NAME/LEDGER_UNDER/TO/BY columns
# (AK:AO) that get populated once a row is actioned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$decisions = @{
    2  = "ACCEPTED"
    3  = "ACCEPTED"
    4  = "ACCEPTED"
    5  = "ACCEPTED"
    6  = "HOLD"
    7  = "HOLD"
    8  = "HOLD"
    9  = "ACCEPTED"
    10 = "ACCEPTED"
    11 = "ACCEPTED"
    12 = "ACCEPTED"
    13 = "ACCEPTED"
    14 = "ACCEPTED"
    15 = "HOLD"
    16 = "ACCEPTED"
    17 = "ACCEPTED"
    18 = "REJECTED"
    19 = "ACCEPTED"
    20 = "HOLD"
    21 = "HOLD"
    22 = "HOLD"
    23 = "HOLD"
    24 = "HOLD"
    25 = "HOLD"
}

foreach ($row in 2..25) {
    $status = $decisions[$row]

    $ws.Range("AI$row").Value = $status
    $ws.Range("AJ$row").Value = $status

    $ws.Range("AK$row").Value = 0
    $ws.Range("AL$row").Value = 0
    $ws.Range("AM$row").Value = 0
    $ws.Range("AN$row").Value = 0
    $ws.Range("AO$row").Value = 0
}
